$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 6185.5454
$ws.Range("I40").Value = 5149.143
$ws.Range("J40").Value = 7999.25
$ws.Range("K40").Value = 5149.143
$ws.Range("L40").Value = 7999.25
$ws.Range("M40").Value = -4974.143
$ws.Range("N40").Value = -8349.25
# Row 51
$ws.Range("H51").Value = 29419700
$ws.Range("I51").Value = 38469536
$ws.Range("J51").Value = 7727.75
$ws.Range("K51").Value = 38469536
$ws.Range("L51").Value = 7727.75
$ws.Range("M51").Value = -38469052
$ws.Range("N51").Value = -8695.75
# Row 64
$ws.Range("H64").Value = 20414136
$ws.Range("I64").Value = 6052.643
$ws.Range("J64").Value = 142862640
$ws.Range("K64").Value = 6052.643
$ws.Range("L64").Value = 142862640
$ws.Range("M64").Value = -5804.643
$ws.Range("N64").Value = -142863136
# Row 67
$ws.Range("H67").Value = 20414136
$ws.Range("I67").Value = 6052.643
$ws.Range("J67").Value = 142862640
$ws.Range("K67").Value = 6052.643
$ws.Range("L67").Value = 142862640
$ws.Range("M67").Value = -5194.643
$ws.Range("N67").Value = -142864356
# Row 74
$ws.Range("H74").Value = 14869.105
$ws.Range("I74").Value = 15139.611
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 15139.611
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -14203.611
$ws.Range("N74").Value = -11872
# Row 77
$ws.Range("H77").Value = 14869.105
$ws.Range("I77").Value = 15139.611
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 75698.05500000001
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -71018.05500000001
$ws.Range("N77").Value = -59360
# Row 106
$ws.Range("H106").Value = 683.1667
$ws.Range("I106").Value = 683.1667
$ws.Range("K106").Value = 683.1667
$ws.Range("M106").Value = -52.16669999999999
# Row 132
$ws.Range("H132").Value = 6569.56
$ws.Range("I132").Value = 6634.9585
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 19904.8755
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -17374.8755
$ws.Range("N132").Value = -20060
# Row 138
$ws.Range("H138").Value = 2717.873
$ws.Range("I138").Value = 2526.2415
$ws.Range("J138").Value = 2881.3235
$ws.Range("K138").Value = 7578.7245
$ws.Range("L138").Value = 8643.970499999999
$ws.Range("M138").Value = -2438.7245
$ws.Range("N138").Value = -18923.9705

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9013
$ws.Range("I32").Value = 9013
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9013
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -8726
# Row 45
$ws.Range("H45").Value = 26103.611
$ws.Range("I45").Value = 34023.77
$ws.Range("K45").Value = 34023.77
$ws.Range("M45").Value = -33646.77
# Row 63
$ws.Range("H63").Value = 4098.75
$ws.Range("J63").Value = 3000
$ws.Range("L63").Value = 3000
$ws.Range("N63").Value = -4372
# Row 66
$ws.Range("H66").Value = 4098.75
$ws.Range("J66").Value = 3000
$ws.Range("L66").Value = 15000
$ws.Range("N66").Value = -21864
# Row 102
$ws.Range("H102").Value = 4963.385
$ws.Range("I102").Value = 4684.091
$ws.Range("K102").Value = 4684.091
$ws.Range("M102").Value = -3062.091

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 29422152
$ws.Range("I20").Value = 35725744
$ws.Range("K20").Value = 35725744
$ws.Range("M20").Value = -35725497
# Row 99
$ws.Range("H99").Value = 4640.9585
$ws.Range("I99").Value = 3164
$ws.Range("K99").Value = 3164
$ws.Range("M99").Value = -1666
# Row 134
$ws.Range("H134").Value = 2407.7856
$ws.Range("I134").Value = 2194.0652
$ws.Range("K134").Value = 6582.1956
$ws.Range("M134").Value = -4047.1956

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2676.25
$ws.Range("I58").Value = 1424.2
$ws.Range("J58").Value = 3570.5715
$ws.Range("K58").Value = 1424.2
$ws.Range("L58").Value = 3570.5715
$ws.Range("M58").Value = -1221.2
$ws.Range("N58").Value = -3976.5715
# Row 86
$ws.Range("H86").Value = 49998.5
$ws.Range("J86").Value = 49998.5
$ws.Range("L86").Value = 49998.5
$ws.Range("N86").Value = -52244.5
# Row 89
$ws.Range("H89").Value = 49998.5
$ws.Range("J89").Value = 49998.5
$ws.Range("L89").Value = 249992.5
$ws.Range("N89").Value = -261224.5
# Row 99
$ws.Range("H99").Value = 4734.1665
$ws.Range("I99").Value = 4781
$ws.Range("K99").Value = 4781
$ws.Range("M99").Value = -3283
# Row 126
$ws.Range("H126").Value = 4734.1665
$ws.Range("I126").Value = 4781
$ws.Range("K126").Value = 14343
$ws.Range("M126").Value = -11873
# Row 132
$ws.Range("H132").Value = 2290.6667
$ws.Range("I132").Value = 1577.125
$ws.Range("K132").Value = 4731.375
$ws.Range("M132").Value = -2201.375
# Row 134
$ws.Range("H134").Value = 2789.4443
$ws.Range("I134").Value = 2547.3333
$ws.Range("K134").Value = 7641.999899999999
$ws.Range("M134").Value = -5106.999899999999
# Row 136
$ws.Range("H136").Value = 2676.25
$ws.Range("I136").Value = 1424.2
$ws.Range("J136").Value = 3570.5715
$ws.Range("K136").Value = 4272.6
$ws.Range("L136").Value = 10711.7145
$ws.Range("M136").Value = -1722.6
$ws.Range("N136").Value = -15811.7145

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 364.2857
$ws.Range("J12").Value = 363.33334
$ws.Range("L12").Value = 1090.00002
$ws.Range("N12").Value = -1436.00002
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0
# Row 61
$ws.Range("H61").Value = 454.22223
$ws.Range("I61").Value = 178.6
$ws.Range("K61").Value = 535.8
$ws.Range("M61").Value = -320.8
# Row 75
$ws.Range("H75").Value = 1961.3334
$ws.Range("J75").Value = 1961.3334
$ws.Range("L75").Value = 5884.0002
$ws.Range("N75").Value = -7880.0002
# Row 78
$ws.Range("H78").Value = 1961.3334
$ws.Range("J78").Value = 1961.3334
$ws.Range("L78").Value = 17652.0006
$ws.Range("N78").Value = -27636.0006
# Row 80
$ws.Range("H80").Value = 1812
$ws.Range("I80").Value = 524
$ws.Range("J80").Value = 3100
$ws.Range("K80").Value = 1572
$ws.Range("L80").Value = 9300
$ws.Range("M80").Value = -636
$ws.Range("N80").Value = -11172
# Row 83
$ws.Range("H83").Value = 1812
$ws.Range("I83").Value = 524
$ws.Range("J83").Value = 3100
$ws.Range("K83").Value = 4716
$ws.Range("L83").Value = 27900
$ws.Range("M83").Value = -36
$ws.Range("N83").Value = -37260
# Row 107
$ws.Range("H107").Value = 336.33334
$ws.Range("J107").Value = 336.33334
$ws.Range("L107").Value = 1009.00002
$ws.Range("N107").Value = -4849.00002
# Row 139
$ws.Range("H139").Value = 2872.8276
$ws.Range("J139").Value = 3099.4092
$ws.Range("L139").Value = 9298.2276
$ws.Range("N139").Value = -19578.2276

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1172.2307
$ws.Range("I22").Value = 1749.8334
$ws.Range("K22").Value = 1749.8334
$ws.Range("M22").Value = -1454.8334
# Row 27
$ws.Range("H27").Value = 1172.2307
$ws.Range("I27").Value = 1749.8334
$ws.Range("K27").Value = 1749.8334
$ws.Range("M27").Value = -1642.8334
# Row 40
$ws.Range("H40").Value = 4267.927
$ws.Range("I40").Value = 4115.353
$ws.Range("K40").Value = 4115.353
$ws.Range("M40").Value = -3979.353
# Row 93
$ws.Range("H93").Value = 563.2
$ws.Range("I93").Value = 563.2
$ws.Range("K93").Value = 563.2
$ws.Range("M93").Value = 684.8
# Row 122
$ws.Range("H122").Value = 4703.9565
$ws.Range("I122").Value = 4463.227
$ws.Range("K122").Value = 13389.681
$ws.Range("M122").Value = -10939.681
# Row 132
$ws.Range("H132").Value = 3324.6511
$ws.Range("I132").Value = 2267.0293
$ws.Range("K132").Value = 6801.0879
$ws.Range("M132").Value = -4271.0879
# Row 136
$ws.Range("H136").Value = 9700.799999999999
$ws.Range("I136").Value = 5286.857
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 15860.571
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -13310.571
$ws.Range("N136").Value = -65100

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 6493.75
$ws.Range("I62").Value = 4737
$ws.Range("K62").Value = 4737
$ws.Range("M62").Value = -4113
# Row 65
$ws.Range("H65").Value = 6493.75
$ws.Range("I65").Value = 4737
$ws.Range("K65").Value = 23685
$ws.Range("M65").Value = -20565
# Row 96
$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("N96").Value = -7746
# Row 107
$ws.Range("H107").Value = 714.26666
$ws.Range("I107").Value = 716.53845
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 2149.61535
$ws.Range("L107").Value = 2098.5
$ws.Range("M107").Value = -229.61535
$ws.Range("N107").Value = -5938.5
# Row 135
$ws.Range("H135").Value = 72727.17999999999
$ws.Range("J135").Value = 72727.17999999999
$ws.Range("L135").Value = 72727.17999999999
$ws.Range("N135").Value = -82867.17999999999
# Row 138
$ws.Range("H138").Value = 124997.5
$ws.Range("J138").Value = 124997.5
$ws.Range("L138").Value = 124997.5
$ws.Range("N138").Value = -135277.5
# Row 140
$ws.Range("H140").Value = 95128.5
$ws.Range("J140").Value = 95128.5
$ws.Range("L140").Value = 95128.5
$ws.Range("N140").Value = -105488.5
# Row 141
$ws.Range("H141").Value = 70141.57000000001
$ws.Range("J141").Value = 70141.57000000001
$ws.Range("L141").Value = 70141.57000000001
$ws.Range("N141").Value = -80501.57000000001
